$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Clear()

$ws.Range("A1").Value = "negative"
$ws.Range("J1").Value = "positive"
$ws.Range("A2").Value = "name"
$ws.Range("B2").Value = "anchor score"
$ws.Range("C2").Value = "type occurences"
$ws.Range("D2").Value = "total occurences"
$ws.Range("E2").Value = "+%"
$ws.Range("F2").Value = "-%"
$ws.Range("G2").Value = "both"
$ws.Range("H2").Value = "normal"
$ws.Range("J2").Value = "name"
$ws.Range("K2").Value = "anchor score"
$ws.Range("L2").Value = "type occurences"
$ws.Range("M2").Value = "total occurences"
$ws.Range("N2").Value = "+%"
$ws.Range("O2").Value = "-%"
$ws.Range("P2").Value = "both"
$ws.Range("Q2").Value = "normal"
$ws.Range("A3").Value = "poorly"
$ws.Range("B3").Value = 0.9565217391304348
$ws.Range("C3").Value = 44
$ws.Range("D3").Value = 44
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = 2
$ws.Range("J3").Value = "awesome"
$ws.Range("K3").Value = 0.8
$ws.Range("L3").Value = 52
$ws.Range("M3").Value = 52
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = 13
$ws.Range("A4").Value = "disappointing"
$ws.Range("B4").Value = 0.8863636363636364
$ws.Range("C4").Value = 39
$ws.Range("D4").Value = 39
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = 5
$ws.Range("J4").Value = "wonderful"
$ws.Range("K4").Value = 0.7857142857142857
$ws.Range("L4").Value = 44
$ws.Range("M4").Value = 44
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = $false
$ws.Range("Q4").Value = 12
$ws.Range("A5").Value = "however"
$ws.Range("B5").Value = 0.765625
$ws.Range("C5").Value = 49
$ws.Range("D5").Value = 49
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = 15
$ws.Range("J5").Value = "favorite"
$ws.Range("K5").Value = 0.6129032258064516
$ws.Range("L5").Value = 57
$ws.Range("M5").Value = 57
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $false
$ws.Range("Q5").Value = 36
$ws.Range("A6").Value = "disappointed"
$ws.Range("B6").Value = 0.7311827956989247
$ws.Range("C6").Value = 136
$ws.Range("D6").Value = 136
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = 50
$ws.Range("J6").Value = "classic"
$ws.Range("K6").Value = 0.6037735849056604
$ws.Range("L6").Value = 32
$ws.Range("M6").Value = 32
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = $false
$ws.Range("Q6").Value = 21
$ws.Range("A7").Value = "poor"
$ws.Range("B7").Value = 0.704225352112676
$ws.Range("C7").Value = 50
$ws.Range("D7").Value = 50
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = 21
$ws.Range("J7").Value = "excellent"
$ws.Range("K7").Value = 0.5
$ws.Range("L7").Value = 32
$ws.Range("M7").Value = 32
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = $false
$ws.Range("Q7").Value = 32
$ws.Range("A8").Value = "broke"
$ws.Range("B8").Value = 0.7038834951456311
$ws.Range("C8").Value = 145
$ws.Range("D8").Value = 145
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = $false
$ws.Range("H8").Value = 61
$ws.Range("J8").Value = "thank"
$ws.Range("K8").Value = 0.4057971014492754
$ws.Range("L8").Value = 28
$ws.Range("M8").Value = 28
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = $false
$ws.Range("Q8").Value = 41
$ws.Range("A9").Value = "waste"
$ws.Range("B9").Value = 0.6621621621621622
$ws.Range("C9").Value = 98
$ws.Range("D9").Value = 98
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = $false
$ws.Range("H9").Value = 50
$ws.Range("J9").Value = "great"
$ws.Range("K9").Value = 0.3352459016393443
$ws.Range("L9").Value = 409
$ws.Range("M9").Value = 409
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = $false
$ws.Range("Q9").Value = 811
$ws.Range("A10").Value = "returned"
$ws.Range("B10").Value = 0.6578947368421053
$ws.Range("C10").Value = 25
$ws.Range("D10").Value = 25
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = $false
$ws.Range("H10").Value = 13
$ws.Range("J10").Value = "love"
$ws.Range("K10").Value = 0.3113342898134864
$ws.Range("L10").Value = 217
$ws.Range("M10").Value = 217
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = $false
$ws.Range("Q10").Value = 480
$ws.Range("A11").Value = "smaller"
$ws.Range("B11").Value = 0.5966386554621849
$ws.Range("C11").Value = 71
$ws.Range("D11").Value = 71
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = $false
$ws.Range("H11").Value = 48
$ws.Range("J11").Value = "loves"
$ws.Range("K11").Value = 0.2551867219917012
$ws.Range("L11").Value = 123
$ws.Range("M11").Value = 123
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = $false
$ws.Range("Q11").Value = 359
$ws.Range("A12").Value = "junk"
$ws.Range("B12").Value = 0.5454545454545454
$ws.Range("C12").Value = 30
$ws.Range("D12").Value = 30
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = $false
$ws.Range("H12").Value = 25
$ws.Range("J12").Value = "best"
$ws.Range("K12").Value = 0.2
$ws.Range("L12").Value = 24
$ws.Range("M12").Value = 24
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = $false
$ws.Range("Q12").Value = 96
$ws.Range("A13").Value = "instead"
$ws.Range("B13").Value = 0.5416666666666666
$ws.Range("C13").Value = 26
$ws.Range("D13").Value = 26
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = $false
$ws.Range("H13").Value = 22
$ws.Range("J13").Value = "perfect"
$ws.Range("K13").Value = 0.1867469879518072
$ws.Range("L13").Value = 31
$ws.Range("M13").Value = 31
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = $false
$ws.Range("Q13").Value = 135
$ws.Range("A14").Value = "missing"
$ws.Range("B14").Value = 0.4888888888888889
$ws.Range("C14").Value = 22
$ws.Range("D14").Value = 22
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = $false
$ws.Range("H14").Value = 23
$ws.Range("J14").Value = "loved"
$ws.Range("K14").Value = 0.1743119266055046
$ws.Range("L14").Value = 57
$ws.Range("M14").Value = 57
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = $false
$ws.Range("Q14").Value = 270
$ws.Range("A15").Value = "small"
$ws.Range("B15").Value = 0.4869565217391305
$ws.Range("C15").Value = 168
$ws.Range("D15").Value = 168
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = $false
$ws.Range("H15").Value = 177
$ws.Range("J15").Value = "friends"
$ws.Range("K15").Value = 0.1164021164021164
$ws.Range("L15").Value = 22
$ws.Range("M15").Value = 22
$ws.Range("N15").Value = 1
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = $false
$ws.Range("Q15").Value = 167
$ws.Range("A16").Value = "paint"
$ws.Range("B16").Value = 0.4761904761904762
$ws.Range("C16").Value = 30
$ws.Range("D16").Value = 30
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = $false
$ws.Range("H16").Value = 33
$ws.Range("J16").Value = "christmas"
$ws.Range("K16").Value = 0.08835341365461848
$ws.Range("L16").Value = 22
$ws.Range("M16").Value = 22
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = $false
$ws.Range("Q16").Value = 227
$ws.Range("A17").Value = "okay"
$ws.Range("B17").Value = 0.4259259259259259
$ws.Range("C17").Value = 23
$ws.Range("D17").Value = 23
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = $false
$ws.Range("H17").Value = 31
$ws.Range("J17").Value = "fun"
$ws.Range("K17").Value = 0.08501314636283962
$ws.Range("L17").Value = 97
$ws.Range("M17").Value = 97
$ws.Range("N17").Value = 1
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = $false
$ws.Range("Q17").Value = 1044
$ws.Range("A18").Value = "guess"
$ws.Range("B18").Value = 0.4259259259259259
$ws.Range("C18").Value = 23
$ws.Range("D18").Value = 23
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = $false
$ws.Range("H18").Value = 31
$ws.Range("J18").Value = "game"
$ws.Range("K18").Value = 0.02790395846852693
$ws.Range("L18").Value = 43
$ws.Range("M18").Value = 43
$ws.Range("N18").Value = 1
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = $false
$ws.Range("Q18").Value = 1498
$ws.Range("A19").Value = "broken"
$ws.Range("B19").Value = 0.4216867469879518
$ws.Range("C19").Value = 35
$ws.Range("D19").Value = 35
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = $false
$ws.Range("H19").Value = 48
$ws.Range("A20").Value = "plastic"
$ws.Range("B20").Value = 0.4173228346456693
$ws.Range("C20").Value = 53
$ws.Range("D20").Value = 53
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = $false
$ws.Range("H20").Value = 74
$ws.Range("A21").Value = "less"
$ws.Range("B21").Value = 0.4
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = $false
$ws.Range("H21").Value = 36
$ws.Range("A22").Value = "di"
$ws.Range("B22").Value = 0.375
$ws.Range("C22").Value = 24
$ws.Range("D22").Value = 24
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = $false
$ws.Range("H22").Value = 40
$ws.Range("A23").Value = "tried"
$ws.Range("B23").Value = 0.360655737704918
$ws.Range("C23").Value = 22
$ws.Range("D23").Value = 22
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = $false
$ws.Range("H23").Value = 39
$ws.Range("A24").Value = "difficult"
$ws.Range("B24").Value = 0.3483146067415731
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = $false
$ws.Range("H24").Value = 58
$ws.Range("A25").Value = "apart"
$ws.Range("B25").Value = 0.3473684210526316
$ws.Range("C25").Value = 33
$ws.Range("D25").Value = 33
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = $false
$ws.Range("H25").Value = 62
$ws.Range("A26").Value = "ok"
$ws.Range("B26").Value = 0.3203125
$ws.Range("C26").Value = 41
$ws.Range("D26").Value = 41
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = $false
$ws.Range("H26").Value = 87
$ws.Range("A27").Value = "cheap"
$ws.Range("B27").Value = 0.2938388625592417
$ws.Range("C27").Value = 62
$ws.Range("D27").Value = 62
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = $false
$ws.Range("H27").Value = 149
$ws.Range("A28").Value = "thought"
$ws.Range("B28").Value = 0.2920792079207921
$ws.Range("C28").Value = 59
$ws.Range("D28").Value = 59
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = $false
$ws.Range("H28").Value = 143
$ws.Range("A29").Value = "size"
$ws.Range("B29").Value = 0.2474226804123711
$ws.Range("C29").Value = 48
$ws.Range("D29").Value = 48
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = $false
$ws.Range("H29").Value = 146
$ws.Range("A30").Value = "though"
$ws.Range("B30").Value = 0.2393162393162393
$ws.Range("C30").Value = 28
$ws.Range("D30").Value = 28
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = $false
$ws.Range("H30").Value = 89
$ws.Range("A31").Value = "bit"
$ws.Range("B31").Value = 0.2244897959183673
$ws.Range("C31").Value = 22
$ws.Range("D31").Value = 22
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = $false
$ws.Range("H31").Value = 76
$ws.Range("A32").Value = "would"
$ws.Range("B32").Value = 0.1899109792284867
$ws.Range("C32").Value = 128
$ws.Range("D32").Value = 128
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 1
$ws.Range("G32").Value = $false
$ws.Range("H32").Value = 546
$ws.Range("A33").Value = "item"
$ws.Range("B33").Value = 0.1884057971014493
$ws.Range("C33").Value = 52
$ws.Range("D33").Value = 52
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = $false
$ws.Range("H33").Value = 224
$ws.Range("A34").Value = "money"
$ws.Range("B34").Value = 0.1708860759493671
$ws.Range("C34").Value = 54
$ws.Range("D34").Value = 54
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 1
$ws.Range("G34").Value = $false
$ws.Range("H34").Value = 262
$ws.Range("A35").Value = "work"
$ws.Range("B35").Value = 0.1708860759493671
$ws.Range("C35").Value = 54
$ws.Range("D35").Value = 54
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 1
$ws.Range("G35").Value = $false
$ws.Range("H35").Value = 262
$ws.Range("A36").Value = "product"
$ws.Range("B36").Value = 0.1409691629955947
$ws.Range("C36").Value = 64
$ws.Range("D36").Value = 64
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 1
$ws.Range("G36").Value = $false
$ws.Range("H36").Value = 390
$ws.Range("A37").Value = "price"
$ws.Range("B37").Value = 0.138328530259366
$ws.Range("C37").Value = 48
$ws.Range("D37").Value = 49
$ws.Range("E37").Value = 0.02
$ws.Range("F37").Value = 0.98
$ws.Range("G37").Value = $true
$ws.Range("H37").Value = 299
$ws.Range("A38").Value = "hard"
$ws.Range("B38").Value = 0.12
$ws.Range("C38").Value = 24
$ws.Range("D38").Value = 24
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = 1
$ws.Range("G38").Value = $false
$ws.Range("H38").Value = 176
$ws.Range("A39").Value = "3"
$ws.Range("B39").Value = 0.09274193548387097
$ws.Range("C39").Value = 23
$ws.Range("D39").Value = 23
$ws.Range("E39").Value = 0
$ws.Range("F39").Value = 1
$ws.Range("G39").Value = $false
$ws.Range("H39").Value = 225
$ws.Range("A40").Value = "2"
$ws.Range("B40").Value = 0.08239700374531835
$ws.Range("C40").Value = 22
$ws.Range("D40").Value = 22
$ws.Range("E40").Value = 0
$ws.Range("F40").Value = 1
$ws.Range("G40").Value = $false
$ws.Range("H40").Value = 245
$ws.Range("A41").Value = "use"
$ws.Range("B41").Value = 0.0684931506849315
$ws.Range("C41").Value = 25
$ws.Range("D41").Value = 25
$ws.Range("E41").Value = 0
$ws.Range("F41").Value = 1
$ws.Range("G41").Value = $false
$ws.Range("H41").Value = 340
$ws.Range("A42").Value = "buy"
$ws.Range("B42").Value = 0.06478873239436619
$ws.Range("C42").Value = 23
$ws.Range("D42").Value = 23
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 1
$ws.Range("G42").Value = $false
$ws.Range("H42").Value = 332
$ws.Range("A43").Value = "much"
$ws.Range("B43").Value = 0.06018518518518518
$ws.Range("C43").Value = 26
$ws.Range("D43").Value = 28
$ws.Range("E43").Value = 0.07000000000000001
$ws.Range("F43").Value = 0.9299999999999999
$ws.Range("G43").Value = $true
$ws.Range("H43").Value = 406
$ws.Range("A44").Value = "like"
$ws.Range("B44").Value = 0.05921052631578947
$ws.Range("C44").Value = 36
$ws.Range("D44").Value = 36
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 1
$ws.Range("G44").Value = $false
$ws.Range("H44").Value = 572
$ws.Range("A45").Value = "little"
$ws.Range("B45").Value = 0.05790645879732739
$ws.Range("C45").Value = 26
$ws.Range("D45").Value = 26
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = 1
$ws.Range("G45").Value = $false
$ws.Range("H45").Value = 423
